$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,13
$arr[0,0] = 1.02
$arr[0,1] = 1.058617298265595
$arr[0,2] = 1.060041922632564
$arr[0,3] = 0.992614727750844
$arr[0,4] = 1.069491942642792
$arr[0,5] = 1
$arr[0,6] = $Null
$arr[0,7] = 1.04683730620434
$arr[0,8] = 1.06360758806204
$arr[0,9] = 1.062769780995098
$arr[0,10] = 0.9955398523335997
$arr[0,11] = 1.072194306124119
$arr[0,12] = 1.065118032834584
$arr[1,0] = 1.02
$arr[1,1] = 1.060249052517938
$arr[1,2] = 1.061291753815073
$arr[1,3] = 0.9936372048519299
$arr[1,4] = 1.070891179892598
$arr[1,5] = 1
$arr[1,6] = $Null
$arr[1,7] = 1.047290171611087
$arr[1,8] = 1.064889004872702
$arr[1,9] = 1.063832859810783
$arr[1,10] = 0.9963617723202687
$arr[1,11] = 1.073408278874931
$arr[1,12] = 1.066401269404097
$arr[2,0] = 1.02
$arr[2,1] = 1.061302404144411
$arr[2,2] = 1.062098094487935
$arr[2,3] = 0.9942998659930998
$arr[2,4] = 1.07179443861311
$arr[2,5] = 1
$arr[2,6] = $Null
$arr[2,7] = 1.047580505005246
$arr[2,8] = 1.065715307693802
$arr[2,9] = 1.064517803293111
$arr[2,10] = 0.9968940712668347
$arr[2,11] = 1.074191130198527
$arr[2,12] = 1.067228745669981
$arr[3,0] = 1.02
$arr[3,1] = 1.061744645708584
$arr[3,2] = 1.062436517494834
$arr[3,3] = 0.994578699834602
$arr[3,4] = 1.072173664696835
$arr[3,5] = 1
$arr[3,6] = $Null
$arr[3,7] = 1.047701918313987
$arr[3,8] = 1.066062009937999
$arr[3,9] = 1.064805057359582
$arr[3,10] = 0.9971179600053012
$arr[3,11] = 1.074519609796732
$arr[3,12] = 1.067575940271132
$arr[4,0] = 1.02
$arr[4,1] = 1.061818865900826
$arr[4,2] = 1.062493307447477
$arr[4,3] = 0.994625531979634
$arr[4,4] = 1.072237309130241
$arr[4,5] = 1
$arr[4,6] = $Null
$arr[4,7] = 1.04772226655338
$arr[4,8] = 1.066120183424506
$arr[4,9] = 1.064853247976187
$arr[4,10] = 0.9971555583673455
$arr[4,11] = 1.07457472619451
$arr[4,12] = 1.067634196370663
$arr[5,0] = 1.02
$arr[5,1] = 1.061308315690073
$arr[5,2] = 1.062102618712206
$arr[5,3] = 0.9943035907978918
$arr[5,4] = 1.071799507818334
$arr[5,5] = 1
$arr[5,6] = $Null
$arr[5,7] = 1.047582129856574
$arr[5,8] = 1.065719942987624
$arr[5,9] = 1.064521644320973
$arr[5,10] = 0.9968970624459044
$arr[5,11] = 1.074195521827441
$arr[5,12] = 1.067233387546452
$arr[6,0] = 1.02
$arr[6,1] = 1.059169282444112
$arr[6,2] = 1.060464807350346
$arr[6,3] = 0.9929600610674297
$arr[6,4] = 1.069965270095852
$arr[6,5] = 1
$arr[6,6] = $Null
$arr[6,7] = 1.046990915663224
$arr[6,8] = 1.064041246469577
$arr[6,9] = 1.063129666697585
$arr[6,10] = 0.9958175282591057
$arr[6,11] = 1.072605131716844
$arr[6,12] = 1.065552307086801
$arr[7,0] = 1.02
$arr[7,1] = 1.055380342339274
$arr[7,2] = 1.057560155665237
$arr[7,3] = 0.9906006454969559
$arr[7,4] = 1.066716313888461
$arr[7,5] = 1
$arr[7,6] = $Null
$arr[7,7] = 1.045928267398266
$arr[7,8] = 1.061060857180596
$arr[7,9] = 1.060653969101317
$arr[7,10] = 0.9939188001724441
$arr[7,11] = 1.06978184263405
$arr[7,12] = 1.062567685303174
$arr[8,0] = 1.02
$arr[8,1] = 1.052840407515341
$arr[8,2] = 1.055610692466955
$arr[8,3] = 0.989033133672735
$arr[8,4] = 1.064538495484569
$arr[8,5] = 1
$arr[8,6] = $Null
$arr[8,7] = 1.045205577076733
$arr[8,8] = 1.059058352233933
$arr[8,9] = 1.05898765710459
$arr[8,10] = 0.9926553831429383
$arr[8,11] = 1.067885135116126
$arr[8,12] = 1.060562336569796
$arr[9,0] = 1.02
$arr[9,1] = 1.051737108344036
$arr[9,2] = 1.054763348559596
$arr[9,3] = 0.988355674866747
$arr[9,4] = 1.063592541333511
$arr[9,5] = 1
$arr[9,6] = $Null
$arr[9,7] = 1.044889208849593
$arr[9,8] = 1.058187420921649
$arr[9,9] = 1.058262259133563
$arr[9,10] = 0.9921088820399291
$arr[9,11] = 1.067060280957539
$arr[9,12] = 1.059690168435152
$arr[10,0] = 1.02
$arr[10,1] = 1.051326755192962
$arr[10,2] = 1.054448114539229
$arr[10,3] = 0.9881042295826724
$arr[10,4] = 1.063240718634691
$arr[10,5] = 1
$arr[10,6] = $Null
$arr[10,7] = 1.044771174309404
$arr[10,8] = 1.057863331060778
$arr[10,9] = 1.057992222630633
$arr[10,10] = 0.9919059725120875
$arr[10,11] = 1.066753347338378
$arr[10,12] = 1.059365618329504
$arr[11,0] = 1.02
$arr[11,1] = 1.051414801945716
$arr[11,2] = 1.054515755766649
$arr[11,3] = 0.9881581567098651
$arr[11,4] = 1.06331620648184
$arr[11,5] = 1
$arr[11,6] = $Null
$arr[11,7] = 1.044796516784496
$arr[11,8] = 1.057932876189837
$arr[11,9] = 1.058050173291017
$arr[11,10] = 0.9919494934313052
$arr[11,11] = 1.066819210515387
$arr[11,12] = 1.059435262220623
$arr[12,0] = 1.02
$arr[12,1] = 1.05170319949733
$arr[12,2] = 1.054737301327487
$arr[12,3] = 0.9883348863814464
$arr[12,4] = 1.063563468867931
$arr[12,5] = 1
$arr[12,6] = $Null
$arr[12,7] = 1.044879462739254
$arr[12,8] = 1.058160643596768
$arr[12,9] = 1.058239949951017
$arr[12,10] = 0.9920921077337197
$arr[12,11] = 1.067034920927136
$arr[12,12] = 1.059663353083398
$arr[13,0] = 1.02
$arr[13,1] = 1.051880818932438
$arr[13,2] = 1.054873737263003
$arr[13,3] = 0.9884438009545853
$arr[13,4] = 1.063715755003793
$arr[13,5] = 1
$arr[13,6] = $Null
$arr[13,7] = 1.044930499261734
$arr[13,8] = 1.058300900465792
$arr[13,9] = 1.058356798966018
$arr[13,10] = 0.9921799884222134
$arr[13,11] = 1.067167754568777
$arr[13,12] = 1.059803809133264
$arr[14,0] = 1.02
$arr[14,1] = 1.05291355508797
$arr[14,2] = 1.055666859338009
$arr[14,3] = 0.9890781214508737
$arr[14,4] = 1.064601212387803
$arr[14,5] = 1
$arr[14,6] = $Null
$arr[14,7] = 1.045226500548542
$arr[14,8] = 1.059116071351845
$arr[14,9] = 1.059035716961326
$arr[14,10] = 0.9926916645766087
$arr[14,11] = 1.067939801950601
$arr[14,12] = 1.060620137655476
$arr[15,0] = 1.02
$arr[15,1] = 1.05356041861164
$arr[15,2] = 1.056163496512173
$arr[15,3] = 0.9894763578477731
$arr[15,4] = 1.065155840488847
$arr[15,5] = 1
$arr[15,6] = $Null
$arr[15,7] = 1.045411250442858
$arr[15,8] = 1.059626372274316
$arr[15,9] = 1.059460540535935
$arr[15,10] = 0.9930127773692701
$arr[15,11] = 1.068423124705412
$arr[15,12] = 1.06113116326379
$arr[16,0] = 1.02
$arr[16,1] = 1.053937387338139
$arr[16,2] = 1.05645286719911
$arr[16,3] = 0.9897087662937551
$arr[16,4] = 1.065479062322438
$arr[16,5] = 1
$arr[16,6] = $Null
$arr[16,7] = 1.045518680480247
$arr[16,8] = 1.059923653062707
$arr[16,9] = 1.059707959806028
$arr[16,10] = 0.9932001317071766
$arr[16,11] = 1.06870469518865
$arr[16,12] = 1.061428866224999
$arr[17,0] = 1.02
$arr[17,1] = 1.054065867400702
$arr[17,2] = 1.056551482950142
$arr[17,3] = 0.9897880325774039
$arr[17,4] = 1.065589224879267
$arr[17,5] = 1
$arr[17,6] = $Null
$arr[17,7] = 1.04555525526806
$arr[17,8] = 1.060024955861081
$arr[17,9] = 1.059792260375802
$arr[17,10] = 0.993264023964098
$arr[17,11] = 1.068800645487208
$arr[17,12] = 1.061530312884966
$arr[18,0] = 1.02
$arr[18,1] = 1.053491051116377
$arr[18,2] = 1.056110244115674
$arr[18,3] = 0.9894336180355766
$arr[18,4] = 1.065096363529605
$arr[18,5] = 1
$arr[18,6] = $Null
$arr[18,7] = 1.045391462845018
$arr[18,8] = 1.059571660074532
$arr[18,9] = 1.059414999649707
$arr[18,10] = 0.9929783193490043
$arr[18,11] = 1.068371304351123
$arr[18,12] = 1.061076373366406
$arr[19,0] = 1.02
$arr[19,1] = 1.051618288593105
$arr[19,2] = 1.054672075365074
$arr[19,3] = 0.9882828385668255
$arr[19,4] = 1.063490668856508
$arr[19,5] = 1
$arr[19,6] = $Null
$arr[19,7] = 1.044855051651162
$arr[19,8] = 1.058093588072794
$arr[19,9] = 1.058184081856101
$arr[19,10] = 0.9920501090198107
$arr[19,11] = 1.066971414742942
$arr[19,12] = 1.059596202332888
$arr[20,0] = 1.02
$arr[20,1] = 1.050437684448142
$arr[20,2] = 1.053764985584887
$arr[20,3] = 0.9875604150241496
$arr[20,4] = 1.062478477306684
$arr[20,5] = 1
$arr[20,6] = $Null
$arr[20,7] = 1.044514770233232
$arr[20,8] = 1.057160862559593
$arr[20,9] = 1.057406728613322
$arr[20,10] = 0.991467000034148
$arr[20,11] = 1.066088084174288
$arr[20,12] = 1.058662152242473
$arr[21,0] = 1.02
$arr[21,1] = 1.051063845970484
$arr[21,2] = 1.0542461251142
$arr[21,3] = 0.9879432794636459
$arr[21,4] = 1.06301531191224
$arr[21,5] = 1
$arr[21,6] = $Null
$arr[21,7] = 1.044695447541299
$arr[21,8] = 1.05765564435424
$arr[21,9] = 1.057819146216395
$arr[21,10] = 0.9917760702887607
$arr[21,11] = 1.066556657775261
$arr[21,12] = 1.059157636684021
$arr[22,0] = 1.02
$arr[22,1] = 1.053522396339228
$arr[22,2] = 1.056134307537439
$arr[22,3] = 0.9894529299347241
$arr[22,4] = 1.065123239482949
$arr[22,5] = 1
$arr[22,6] = $Null
$arr[22,7] = 1.045400405032917
$arr[22,8] = 1.059596383303074
$arr[22,9] = 1.059435578768893
$arr[22,10] = 0.9929938892766438
$arr[22,11] = 1.068394720800091
$arr[22,12] = 1.061101131704768
$arr[23,0] = 1.019999999999999
$arr[23,1] = 1.05636228292944
$arr[23,2] = 1.058313336410784
$arr[23,3] = 0.9912096547607046
$arr[23,4] = 1.067558294092036
$arr[23,5] = 1
$arr[23,6] = $Null
$arr[23,7] = 1.046205482963907
$arr[23,8] = 1.061834062410455
$arr[23,9] = 1.061296753885314
$arr[23,10] = 0.9944092447426411
$arr[23,11] = 1.070514250599938
$arr[23,12] = 1.063341988573147

$ws.Range("B2:N25").Value = $arr
